$d = $word.ActiveDocument

# 1. Change the second heading from "III: Theoretical Questions" to
#    "III: Normal Environment Evaluator"
$d.Content.Find.Execute("III: Theoretical Questions", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "III: Normal Environment Evaluator", 2)

# 2. Add new text to the last (empty) paragraph and change its indentation
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Text = "3.2+3.3) We implemented the bonus. We did not change the type Closure at all."
$lastPara.Format.LeftIndent = 0
